# Add team record (Wins/Losses/Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy the existing header formatting (bold + border + centered)
# from column AC onto the three new header cells, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-45): every player row gets the same team record values.
for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 30).Value = 83
    $ws.Cells.Item($row, 31).Value = 79
    $ws.Cells.Item($row, 32).Value = 0
}
